# Generate Report for Archive
# Update localization status for the two files that are still being
# translated ("25b6cd32-...md" and "63c391f1-...md") from
# "Ready for handoff" to "In Translation" across the Overview summary
# sheet and each per-language detail sheet. The "8df1a12c-...md" row
# (row 5) keeps its "Ready for handoff" status and is left untouched.

$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn (col E) and de-de (col F) status columns
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"   # 25b6cd32...md, zh-cn
$overview.Range("F3").Value = "In Translation"   # 25b6cd32...md, de-de
$overview.Range("E4").Value = "In Translation"   # 63c391f1...md, zh-cn
$overview.Range("F4").Value = "In Translation"   # 63c391f1...md, de-de

# zh-cn detail sheet: Status column is C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"       # 25b6cd32...md
$zhcn.Range("C4").Value = "In Translation"       # 63c391f1...md

# de-de detail sheet: Status column is C
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"       # 25b6cd32...md
$dede.Range("C4").Value = "In Translation"       # 63c391f1...md
